$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.70%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'27.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.53%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.207"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.13%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.05920"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.74%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'6.676"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.54%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.8664"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.87%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'1.013"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'14.70%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1415"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.98%"
$ws.Range("E9").Style = "Normal"

$ws.Range("B10").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.03610"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'8.55%"
$ws.Range("E10").Style = "Normal"

$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.07190"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.58%"
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.03194"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.68%"
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").Value = "'BitMartToken"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.09229"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.12%"
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = "'BitForexToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.001537"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.11%"
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = "'One"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.0006052"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.58%"
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "'TigerCash"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.005930"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.32%"
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = "'LEO"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.478"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.37%"
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = "'GateToken"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'3.266"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.77%"
$ws.Range("E18").Style = "Normal"

$ws.Range("B19").Value = "'BTSEToken"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'2.224"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.63%"
$ws.Range("E19").Style = "Normal"

$ws.Range("B20").Value = "'BitpandaEcosystemToken"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.3150"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.55%"
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'0.05%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'3.525"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.09%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04163"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.84%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.1400"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.56%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.001217"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.28%"
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'8.76%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0001200"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.06%"
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'33.84%"
$ws.Range("E28").Style = "Normal"

$ws.Range("D40").Value = "'0.03828"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.73%"
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'BKEXToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.1104"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.42%"
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'KickToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.003975"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-23.11%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002460"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.14%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01082"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'14.03%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005431"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.91%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'3.97%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.002174"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-4.16%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.07%"
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'0.07%"
$ws.Range("E50").Style = "Normal"
